# "Final slide, and paper"
#
# Slide 13 ("Author based statistics") — tidy up the bullet-point text box:
#   1) The first bullet's three runs ("Steady amount of research per " /
#      "author; " / "good as number of authors increasing") collapse into a
#      single run once the text is re-typed as one continuous string.
#   2) "Lotka" becomes "Lotka's" (curly apostrophe), keeping the original
#      run's spell-check "err" formatting flag.
#   3) The following run (" Law of 60%") is split into a plain space run and
#      a "Law of 60%" run.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(13)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# 1) Merge "Steady amount of research per " + "author; " + "good as number
#    of authors increasing" into a single run of text.
$bullet1 = $tr.Find("Steady amount of research per author; good as number of authors increasing")
if ($null -eq $bullet1) {
    $bullet1 = $tr.Find("Steady amount of research per")
}
$bullet1.Text = "Steady amount of research per author; good as number of authors increasing"

# 2) "Lotka" -> "Lotka's" on the existing run (preserves its rPr, incl. err="1").
$lotka = $tr.Find("Lotka")
$lotka.Text = "Lotka" + [char]0x2019 + "s"

# 3) Re-type the tail so it splits into " " and "Law of 60%" runs.
$law = $tr.Find("Law of 60%")
$law.Text = "Law of 60%"
